# Generate Report for Handoff
# Re-running the handoff report unified the "Ready for handoff" rows'
# latest handoff/handback timestamps to the freshest run time.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D (Latest Handoff Date) ---
$ws = $wb.Worksheets.Item("Overview")
foreach ($r in @(7, 10, 11, 12, 13, 14, 15, 16)) {
    $ws.Cells.Item($r, 4).Value = "2016-03-22 22:35:59"
}

# --- zh-cn sheet: column E (Latest Handoff Datetime) ---
$ws = $wb.Worksheets.Item("zh-cn")
foreach ($r in @(7, 10, 11, 12, 13, 14, 15, 16)) {
    $ws.Cells.Item($r, 5).Value = "2016-03-22 22:35:55"
}

# --- de-de sheet: column E (Latest Handoff Datetime) ---
$ws = $wb.Worksheets.Item("de-de")
foreach ($r in @(7, 10, 11, 12, 13, 14, 15, 16)) {
    $ws.Cells.Item($r, 5).Value = "2016-03-22 22:35:59"
}
